$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '310.23'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.58%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '10'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '37.12'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-3.15%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '10'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.124'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.27%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '10'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07757'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.96%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '10'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.380'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.89%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '10'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.209'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.46%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '10'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.879'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-7.61%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '10'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.946'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-3.90%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '10'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9191'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.33%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '10'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1211'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-5.42%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '10'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1896'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.61%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '10'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09194'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '4.78%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '10'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03433'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.19%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '10'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09690'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.42%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '10'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001367'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.12%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '10'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005929'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-5.96%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '10'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.559'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.05%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '10'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3407'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.35%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '10'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.289'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '5.15%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '10'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.54%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '10'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2593'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2.56%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '10'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.02107'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '5,596.73%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '10'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.04359'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.47%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '10'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001198'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-1.95%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '10'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004246'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-8.42%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '10'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001302'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-63.75%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '10'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '10'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '10'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '10'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '10'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '10'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '10'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '10'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '10'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '10'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '10'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '10'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02086'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-5.92%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '10'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05027'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.01%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '10'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007682'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.38%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '10'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-1.27%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '10'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1345'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.95%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '10'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002084'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '2.50%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '10'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009598'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '8.50%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '10'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006705'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.37%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '10'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.45%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '10'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-0.32%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '10'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002939'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-2.51%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '10'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.45%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '10'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.45%'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '10'
